# Apply updated TPM-derived values to rows 2-17 (columns G:T) per the new NATMI output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.713280333333334
$ws.Range("H2").Value = 14.139841
$ws.Range("I2").Value = 0.6942627767023021
$ws.Range("J2").Value = 0.6942627767023022
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 17.88780538591745
$ws.Range("R2").Value = 160.990248473257
$ws.Range("S2").Value = 0.007249587415118269
$ws.Range("T2").Value = 0.00724958741511827

# Row 3
$ws.Range("G3").Value = 4.713280333333334
$ws.Range("H3").Value = 14.139841
$ws.Range("I3").Value = 0.6942627767023021
$ws.Range("J3").Value = 0.6942627767023022
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 1147.101193433093
$ws.Range("R3").Value = 10323.91074089784
$ws.Range("S3").Value = 0.464898303417738
$ws.Range("T3").Value = 0.4648983034177381

# Row 4
$ws.Range("G4").Value = 4.713280333333334
$ws.Range("H4").Value = 14.139841
$ws.Range("I4").Value = 0.6942627767023021
$ws.Range("J4").Value = 0.6942627767023022
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 140.4622158406703
$ws.Range("R4").Value = 1264.159942566033
$ws.Range("S4").Value = 0.0569266654175376
$ws.Range("T4").Value = 0.0569266654175376

# Row 5
$ws.Range("G5").Value = 4.713280333333334
$ws.Range("H5").Value = 14.139841
$ws.Range("I5").Value = 0.6942627767023021
$ws.Range("J5").Value = 0.6942627767023022
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 407.5893661655442
$ws.Range("R5").Value = 3668.304295489898
$ws.Range("S5").Value = 0.1651882204519082
$ws.Range("T5").Value = 0.1651882204519082

# Row 6
$ws.Range("I6").Value = 0.1843120478188439
$ws.Range("J6").Value = 0.1843120478188439
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 4.748833082084001
$ws.Range("R6").Value = 42.73949773875601
$ws.Range("S6").Value = 0.001924611756760107
$ws.Range("T6").Value = 0.001924611756760107

# Row 7
$ws.Range("I7").Value = 0.1843120478188439
$ws.Range("J7").Value = 0.1843120478188439
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("Q7").Value = 304.5310466180321
$ws.Range("R7").Value = 2740.779419562289
$ws.Range("S7").Value = 0.1234206430271741
$ws.Range("T7").Value = 0.1234206430271741

# Row 8
$ws.Range("I8").Value = 0.1843120478188439
$ws.Range("J8").Value = 0.1843120478188439
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 37.289740299396
$ws.Range("R8").Value = 335.607662694564
$ws.Range("S8").Value = 0.01511282273902403
$ws.Range("T8").Value = 0.01511282273902403

# Row 9
$ws.Range("I9").Value = 0.1843120478188439
$ws.Range("J9").Value = 0.1843120478188439
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 108.206335220776
$ws.Range("R9").Value = 973.8570169869841
$ws.Range("S9").Value = 0.04385397029588561
$ws.Range("T9").Value = 0.04385397029588562

# Row 10
$ws.Range("G10").Value = 0.6263116666666667
$ws.Range("H10").Value = 1.878935
$ws.Range("I10").Value = 0.09225525452111802
$ws.Range("J10").Value = 0.09225525452111803
$ws.Range("M10").Value = 3.795192333333334
$ws.Range("N10").Value = 11.385577
$ws.Range("O10").Value = 0.01044213755712683
$ws.Range("P10").Value = 0.01044213755712683
$ws.Range("Q10").Value = 2.376973235610556
$ws.Range("R10").Value = 21.392759120495
$ws.Range("S10").Value = 0.0009633420580772615
$ws.Range("T10").Value = 0.0009633420580772617

# Row 11
$ws.Range("G11").Value = 0.6263116666666667
$ws.Range("H11").Value = 1.878935
$ws.Range("I11").Value = 0.09225525452111802
$ws.Range("J11").Value = 0.09225525452111803
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.6696287328350964
$ws.Range("P11").Value = 0.6696287328350964
$ws.Range("Q11").Value = 152.4294778762511
$ws.Range("R11").Value = 1371.86530088626
$ws.Range("S11").Value = 0.06177676918235556
$ws.Range("T11").Value = 0.06177676918235556

# Row 12
$ws.Range("G12").Value = 0.6263116666666667
$ws.Range("H12").Value = 1.878935
$ws.Range("I12").Value = 0.09225525452111802
$ws.Range("J12").Value = 0.09225525452111803
$ws.Range("M12").Value = 29.801371
$ws.Range("N12").Value = 89.404113
$ws.Range("O12").Value = 0.08199584844219236
$ws.Range("P12").Value = 0.08199584844219235
$ws.Range("Q12").Value = 18.66494633996167
$ws.Range("R12").Value = 167.984517059655
$ws.Range("S12").Value = 0.007564547867709475
$ws.Range("T12").Value = 0.007564547867709475

# Row 13
$ws.Range("G13").Value = 0.6263116666666667
$ws.Range("H13").Value = 1.878935
$ws.Range("I13").Value = 0.09225525452111802
$ws.Range("J13").Value = 0.09225525452111803
$ws.Range("M13").Value = 86.47679266666667
$ws.Range("N13").Value = 259.430378
$ws.Range("O13").Value = 0.2379332811655844
$ws.Range("P13").Value = 0.2379332811655844
$ws.Range("Q13").Value = 54.16142414304778
$ws.Range("R13").Value = 487.45281728743
$ws.Range("S13").Value = 0.02195059541297572
$ws.Range("T13").Value = 0.02195059541297573

# Row 14
$ws.Range("G14").Value = 0.1980316666666667
$ws.Range("H14").Value = 0.594095
$ws.Range("I14").Value = 0.02916992095773596
$ws.Range("J14").Value = 0.02916992095773596
$ws.Range("M14").Value = 3.795192333333334
$ws.Range("N14").Value = 11.385577
$ws.Range("O14").Value = 0.01044213755712683
$ws.Range("P14").Value = 0.01044213755712683
$ws.Range("Q14").Value = 0.7515682630905558
$ws.Range("R14").Value = 6.764114367815001
$ws.Range("S14").Value = 0.0003045963271711958
$ws.Range("T14").Value = 0.0003045963271711958

# Row 15
$ws.Range("G15").Value = 0.1980316666666667
$ws.Range("H15").Value = 0.594095
$ws.Range("I15").Value = 0.02916992095773596
$ws.Range("J15").Value = 0.02916992095773596
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.6696287328350964
$ws.Range("P15").Value = 0.6696287328350964
$ws.Range("Q15").Value = 48.19623385529112
$ws.Range("R15").Value = 433.7661046976201
$ws.Range("S15").Value = 0.01953301720782865
$ws.Range("T15").Value = 0.01953301720782865

# Row 16
$ws.Range("G16").Value = 0.1980316666666667
$ws.Range("H16").Value = 0.594095
$ws.Range("I16").Value = 0.02916992095773596
$ws.Range("J16").Value = 0.02916992095773596
$ws.Range("M16").Value = 29.801371
$ws.Range("N16").Value = 89.404113
$ws.Range("O16").Value = 0.08199584844219236
$ws.Range("P16").Value = 0.08199584844219235
$ws.Range("Q16").Value = 5.901615168081667
$ws.Range("R16").Value = 53.114536512735
$ws.Range("S16").Value = 0.002391812417921248
$ws.Range("T16").Value = 0.002391812417921248

# Row 17
$ws.Range("G17").Value = 0.1980316666666667
$ws.Range("H17").Value = 0.594095
$ws.Range("I17").Value = 0.02916992095773596
$ws.Range("J17").Value = 0.02916992095773596
$ws.Range("M17").Value = 86.47679266666667
$ws.Range("N17").Value = 259.430378
$ws.Range("O17").Value = 0.2379332811655844
$ws.Range("P17").Value = 0.2379332811655844
$ws.Range("Q17").Value = 17.12514337976778
$ws.Range("R17").Value = 154.12629041791
$ws.Range("S17").Value = 0.006940495004814862
$ws.Range("T17").Value = 0.006940495004814863

